# "Version 2." -> "Version 1."
#
# Original runs in the paragraph:
#   [proofErr spellStart] "Versi" | "on" [proofErr spellEnd] " 2" [bookmark _GoBack] "."
#
# Target runs:
#   [proofErr spellStart] "Version" [proofErr spellEnd] " 1." [bookmark _GoBack]
#
# i.e. "Versi"+"on" merge into a single "Version" run, and " 2"+"."
# merge into a single " 1." run with the (hidden) _GoBack bookmark ending up
# right after it instead of in between.

$d = $word.ActiveDocument

# --- Merge the "Versi" / "on" runs into a single "Version" run ---
# Replacing with exactly the same text wouldn't register as a change (no
# run merge would happen), so nudge it through a temporary value first.
$rVersion = $d.Range(0, 7)
$rVersion.Text = "Versionx"
$rVersionFix = $d.Range(0, 8)
$rVersionFix.Text = "Version"

# --- Drop the trailing "." run first ---
# It sits right after the bookmark, so deleting it on its own (without the
# delete crossing the bookmark) leaves the bookmark anchored at the end of
# the paragraph instead of destroying it.
$rDot = $d.Range(9, 10)
$rDot.Delete()

# --- Turn " 2" into " 1." ---
# This run ends right where the bookmark now sits, so rewriting it doesn't
# cross the bookmark either - the bookmark simply stays put right after the
# newly written text, i.e. right after the merged " 1." run.
$rNum = $d.Range(7, 9)
$rNum.Text = " 1."
